# Daily attendance processing - 2025-11-16 18:50:30
# Reorders the "Recorded By" (column G) comma-separated list so that entries
# starting with "System, ..." have their items reversed (System moved to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.StartsWith("System,")) {
        $parts = $val -split ", "
        $reversedParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversedParts)
    }
}
